$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 41.519606
$ws.Range("H2").Value = 124.558818
$ws.Range("I2").Value = 0.7305114279806179
$ws.Range("J2").Value = 0.7630546295388222
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 34.52018166666667
$ws.Range("N2").Value = 103.560545
$ws.Range("O2").Value = 0.7811694070883366
$ws.Range("P2").Value = 0.8076479349223241
$ws.Range("Q2").Value = 1433.264341848423
$ws.Range("R2").Value = 12899.37907663581
$ws.Range("S2").Value = 0.5706531790668734
$ws.Range("T2").Value = 0.6162794957799488
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 41.519606
$ws.Range("H3").Value = 124.558818
$ws.Range("I3").Value = 0.7305114279806179
$ws.Range("J3").Value = 0.7630546295388222
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.323898
$ws.Range("N3").Value = 15.971694
$ws.Range("O3").Value = 0.1204763718863815
$ws.Range("P3").Value = 0.1245600404701547
$ws.Range("Q3").Value = 221.046147344188
$ws.Range("R3").Value = 1989.415326097692
$ws.Range("S3").Value = 0.0880093664646445
$ws.Range("T3").Value = 0.09504611553629458
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 41.519606
$ws.Range("H4").Value = 124.558818
$ws.Range("I4").Value = 0.7305114279806179
$ws.Range("J4").Value = 0.7630546295388222
$ws.Range("K4").Value = 2
$ws.Range("M4").Value = 4.3463115
$ws.Range("N4").Value = 8.692622999999999
$ws.Range("O4").Value = 0.09835422102528193
$ws.Range("P4").Value = 0.06779202460752111
$ws.Range("Q4").Value = 180.457141033269
$ws.Range("R4").Value = 1082.742846199614
$ws.Range("S4").Value = 0.07184888244910001
$ws.Range("T4").Value = 0.05172901822257873
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 7.183879666666667
$ws.Range("H5").Value = 21.551639
$ws.Range("I5").Value = 0.1263958572665066
$ws.Range("J5").Value = 0.1320266053993819
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 34.52018166666667
$ws.Range("N5").Value = 103.560545
$ws.Range("O5").Value = 0.7811694070883366
$ws.Range("P5").Value = 0.8076479349223241
$ws.Range("Q5").Value = 247.9888311648061
$ws.Range("R5").Value = 2231.899480483255
$ws.Range("S5").Value = 0.098736576879299
$ws.Range("T5").Value = 0.1066310152056154
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 7.183879666666667
$ws.Range("H6").Value = 21.551639
$ws.Range("I6").Value = 0.1263958572665066
$ws.Range("J6").Value = 0.1320266053993819
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.323898
$ws.Range("N6").Value = 15.971694
$ws.Range("O6").Value = 0.1204763718863815
$ws.Range("P6").Value = 0.1245600404701547
$ws.Range("Q6").Value = 38.24624258960733
$ws.Range("R6").Value = 344.216183306466
$ws.Range("S6").Value = 0.01522771430493764
$ws.Range("T6").Value = 0.01644523931168415
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 7.183879666666667
$ws.Range("H7").Value = 21.551639
$ws.Range("I7").Value = 0.1263958572665066
$ws.Range("J7").Value = 0.1320266053993819
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 4.3463115
$ws.Range("N7").Value = 8.692622999999999
$ws.Range("O7").Value = 0.09835422102528193
$ws.Range("P7").Value = 0.06779202460752111
$ws.Range("Q7").Value = 31.2233788098495
$ws.Range("R7").Value = 187.340272859097
$ws.Range("S7").Value = 0.01243156608226998
$ws.Range("T7").Value = 0.008950350882082379
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.229894
$ws.Range("H8").Value = 0.6896819999999999
$ws.Range("I8").Value = 0.004044840748830231
$ws.Range("J8").Value = 0.004225032410066655
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 34.52018166666667
$ws.Range("N8").Value = 103.560545
$ws.Range("O8").Value = 0.7811694070883366
$ws.Range("P8").Value = 0.8076479349223241
$ws.Range("Q8").Value = 7.935982644076665
$ws.Range("R8").Value = 71.42384379668998
$ws.Range("S8").Value = 0.003159705849530455
$ws.Range("T8").Value = 0.003412338700970224
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.229894
$ws.Range("H9").Value = 0.6896819999999999
$ws.Range("I9").Value = 0.004044840748830231
$ws.Range("J9").Value = 0.004225032410066655
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.323898
$ws.Range("N9").Value = 15.971694
$ws.Range("O9").Value = 0.1204763718863815
$ws.Range("P9").Value = 0.1245600404701547
$ws.Range("Q9").Value = 1.223932206812
$ws.Range("R9").Value = 11.015389861308
$ws.Range("S9").Value = 0.0004873077382772605
$ws.Range("T9").Value = 0.0005262702079856177
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.229894
$ws.Range("H10").Value = 0.6896819999999999
$ws.Range("I10").Value = 0.004044840748830231
$ws.Range("J10").Value = 0.004225032410066655
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 4.3463115
$ws.Range("N10").Value = 8.692622999999999
$ws.Range("O10").Value = 0.09835422102528193
$ws.Range("P10").Value = 0.06779202460752111
$ws.Range("Q10").Value = 0.9991909359809997
$ws.Range("R10").Value = 5.995145615885999
$ws.Range("S10").Value = 0.0003978271610225154
$ws.Range("T10").Value = 0.0002864235011108129
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.6310036666666666
$ws.Range("H11").Value = 1.893011
$ws.Range("I11").Value = 0.01110211377241086
$ws.Range("J11").Value = 0.011596696488545
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 34.52018166666667
$ws.Range("N11").Value = 103.560545
$ws.Range("O11").Value = 0.7811694070883366
$ws.Range("P11").Value = 0.8076479349223241
$ws.Range("Q11").Value = 21.78236120566611
$ws.Range("R11").Value = 196.041250850995
$ws.Range("S11").Value = 0.008672631633021444
$ws.Range("T11").Value = 0.009366047970894334
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.6310036666666666
$ws.Range("H12").Value = 1.893011
$ws.Range("I12").Value = 0.01110211377241086
$ws.Range("J12").Value = 0.011596696488545
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 5.323898
$ws.Range("N12").Value = 15.971694
$ws.Range("O12").Value = 0.1204763718863815
$ws.Range("P12").Value = 0.1245600404701547
$ws.Range("Q12").Value = 3.359399158959333
$ws.Range("R12").Value = 30.234592430634
$ws.Range("S12").Value = 0.001337542387569888
$ws.Range("T12").Value = 0.001444484983933265
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.6310036666666666
$ws.Range("H13").Value = 1.893011
$ws.Range("I13").Value = 0.01110211377241086
$ws.Range("J13").Value = 0.011596696488545
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 4.3463115
$ws.Range("N13").Value = 8.692622999999999
$ws.Range("O13").Value = 0.09835422102528193
$ws.Range("P13").Value = 0.06779202460752111
$ws.Range("Q13").Value = 2.7425384929755
$ws.Range("R13").Value = 16.455230957853
$ws.Range("S13").Value = 0.001091939751819524
$ws.Range("T13").Value = 0.000786163533717396
$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 7.2719705
$ws.Range("H14").Value = 14.543941
$ws.Range("I14").Value = 0.1279457602316344
$ws.Range("J14").Value = 0.08909703616318426
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 34.52018166666667
$ws.Range("N14").Value = 103.560545
$ws.Range("O14").Value = 0.7811694070883366
$ws.Range("P14").Value = 0.8076479349223241
$ws.Range("Q14").Value = 251.0297427346408
$ws.Range("R14").Value = 1506.178456407845
$ws.Range("S14").Value = 0.09994731365961229
$ws.Range("T14").Value = 0.0719590372648954
$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 7.2719705
$ws.Range("H15").Value = 14.543941
$ws.Range("I15").Value = 0.1279457602316344
$ws.Range("J15").Value = 0.08909703616318426
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 5.323898
$ws.Range("N15").Value = 15.971694
$ws.Range("O15").Value = 0.1204763718863815
$ws.Range("P15").Value = 0.1245600404701547
$ws.Range("Q15").Value = 38.715229201009
$ws.Range("R15").Value = 232.291375206054
$ws.Range("S15").Value = 0.01541444099095218
$ws.Range("T15").Value = 0.01109793043025707
$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 7.2719705
$ws.Range("H16").Value = 14.543941
$ws.Range("I16").Value = 0.1279457602316344
$ws.Range("J16").Value = 0.08909703616318426
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 4.3463115
$ws.Range("N16").Value = 8.692622999999999
$ws.Range("O16").Value = 0.09835422102528193
$ws.Range("P16").Value = 0.06779202460752111
$ws.Range("Q16").Value = 31.60624901181075
$ws.Range("R16").Value = 126.424996047243
$ws.Range("S16").Value = 0.01258400558106989
$ws.Range("T16").Value = 0.006040068468031786
